# Updated team-specific transition-probability matrix ("Fresno Pacific_A")
# after adding more simulated games, speeding up the simulate-game logic,
# and drafting new optimization logic. The recomputed probabilities for the
# affected states are written back into the corresponding matrix cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3333333333333333
$ws.Range("C2").Value = 0.3333333333333333
$ws.Range("P2").Value = 0.3333333333333333
$ws.Range("P4").Value = 0.5
$ws.Range("S4").Value = 0.5
$ws.Range("J6").Value = 0.4375
$ws.Range("Q6").Value = 0.1875
$ws.Range("R6").Value = 0.0625
$ws.Range("S6").Value = 0.3125
$ws.Range("B7").Value = 0.06666666666666667
$ws.Range("F7").Value = 0.06666666666666667
$ws.Range("J7").Value = 0.06666666666666667
$ws.Range("Q7").Value = 0.06666666666666667
$ws.Range("R7").Value = 0.06666666666666667
$ws.Range("S7").Value = 0.6666666666666666
$ws.Range("B8").Value = 0.04761904761904762
$ws.Range("F8").Value = 0.07142857142857142
$ws.Range("J8").Value = 0.07142857142857142
$ws.Range("O8").Value = 0.02380952380952381
$ws.Range("Q8").Value = 0.1428571428571428
$ws.Range("R8").Value = 0.2142857142857143
$ws.Range("S8").Value = 0.4285714285714285
$ws.Range("B9").Value = 0.08333333333333333
$ws.Range("J9").Value = 0.25
$ws.Range("Q9").Value = 0.08333333333333333
$ws.Range("R9").Value = 0.08333333333333333
$ws.Range("B10").Value = 0.04938271604938271
$ws.Range("D10").Value = 0.02469135802469136
$ws.Range("F10").Value = 0.08641975308641975
$ws.Range("J10").Value = 0.1111111111111111
$ws.Range("O10").Value = 0.02469135802469136
$ws.Range("Q10").Value = 0.1358024691358025
$ws.Range("R10").Value = 0.09876543209876543
$ws.Range("S10").Value = 0.4691358024691358
$ws.Range("G11").Value = 0.1428571428571428
$ws.Range("J11").Value = 0.04761904761904762
$ws.Range("K11").Value = 0.1904761904761905
$ws.Range("L11").Value = 0.6190476190476191
$ws.Range("G12").Value = 0.6428571428571429
$ws.Range("J12").Value = 0.2857142857142857
$ws.Range("L12").Value = 0.07142857142857142
$ws.Range("G13").Value = 0.5714285714285714
$ws.Range("J13").Value = 0.4285714285714285
$ws.Range("H15").Value = 0.3888888888888889
$ws.Range("J15").Value = 0.1666666666666667
$ws.Range("K15").Value = 0.05555555555555555
$ws.Range("M15").Value = 0.05555555555555555
$ws.Range("S15").Value = 0.3333333333333333
$ws.Range("H16").Value = 0.1666666666666667
$ws.Range("J16").Value = 0.5
$ws.Range("K16").Value = 0.1666666666666667
$ws.Range("S16").Value = 0.1666666666666667
$ws.Range("F17").Value = 0.04545454545454546
$ws.Range("H17").Value = 0.2727272727272727
$ws.Range("I17").Value = 0.09090909090909091
$ws.Range("J17").Value = 0.3636363636363636
$ws.Range("K17").Value = 0.09090909090909091
$ws.Range("O17").Value = 0.04545454545454546
$ws.Range("S17").Value = 0.09090909090909091
$ws.Range("F18").Value = 0.04761904761904762
$ws.Range("H18").Value = 0.1428571428571428
$ws.Range("I18").Value = 0.1428571428571428
$ws.Range("J18").Value = 0.3809523809523809
$ws.Range("K18").Value = 0.09523809523809523
$ws.Range("M18").Value = 0.09523809523809523
$ws.Range("O18").Value = 0.04761904761904762
$ws.Range("S18").Value = 0.04761904761904762
$ws.Range("F19").Value = 0.0303030303030303
$ws.Range("H19").Value = 0.2626262626262627
$ws.Range("I19").Value = 0.08080808080808081
$ws.Range("J19").Value = 0.2828282828282828
$ws.Range("K19").Value = 0.1111111111111111
$ws.Range("M19").Value = 0.04040404040404041
$ws.Range("O19").Value = 0.09090909090909091
$ws.Range("S19").Value = 0.101010101010101

